$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15454.546
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9532
$ws.Range("H23").Value = 15454.546
$ws.Range("I23").Value = 10000
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = -9766
$ws.Range("H107").Value = 3631.6667
$ws.Range("I107").Value = 3958
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 3958
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -2038
$ws.Range("N107").Value = -5840
$ws.Range("H113").Value = 11333.333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 11333.333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 11333.333
$ws.Range("N113").Value = -17841.333
$ws.Range("H116").Value = 723014.5600000001
$ws.Range("I116").Value = 2001480.8
$ws.Range("J116").Value = 12755.556
$ws.Range("K116").Value = 2001480.8
$ws.Range("L116").Value = 12755.556
$ws.Range("M116").Value = -1998038.8
$ws.Range("N116").Value = -19639.556
$ws.Range("H125").Value = 2268.7144
$ws.Range("I125").Value = 2548.3333
$ws.Range("J125").Value = 2059
$ws.Range("K125").Value = 22934.9997
$ws.Range("L125").Value = 18531
$ws.Range("M125").Value = -20474.9997
$ws.Range("N125").Value = -23451
$ws.Range("H129").Value = 875.89
$ws.Range("I129").Value = 376.66666
$ws.Range("J129").Value = 891.3299
$ws.Range("K129").Value = 1129.99998
$ws.Range("L129").Value = 2673.9897
$ws.Range("M129").Value = 3870.00002
$ws.Range("N129").Value = -12673.9897
$ws.Range("H137").Value = 1765973.5
$ws.Range("I137").Value = 2507489
$ws.Range("J137").Value = 4874.125
$ws.Range("K137").Value = 7522467
$ws.Range("L137").Value = 14622.375
$ws.Range("M137").Value = -7519917
$ws.Range("N137").Value = -19722.375
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 583.3333
$ws.Range("I2").Value = 583.3333
$ws.Range("K2").Value = 583.3333
$ws.Range("M2").Value = -470.3333
$ws.Range("H32").Value = 3816.209
$ws.Range("I32").Value = 3178.966
$ws.Range("J32").Value = 8515.875
$ws.Range("K32").Value = 3178.966
$ws.Range("L32").Value = 8515.875
$ws.Range("M32").Value = -2891.966
$ws.Range("N32").Value = -9089.875
$ws.Range("H39").Value = 4500
$ws.Range("I39").Value = 4500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 4500
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -3980
$ws.Range("H45").Value = 2651.5715
$ws.Range("I45").Value = 3695.25
$ws.Range("J45").Value = 1260
$ws.Range("K45").Value = 3695.25
$ws.Range("L45").Value = 1260
$ws.Range("M45").Value = -3318.25
$ws.Range("N45").Value = -2014
$ws.Range("H116").Value = 583.3333
$ws.Range("I116").Value = 583.3333
$ws.Range("K116").Value = 583.3333
$ws.Range("M116").Value = 1710.6667
$ws.Range("H132").Value = 1939.5652
$ws.Range("I132").Value = 858.05884
$ws.Range("J132").Value = 5003.8335
$ws.Range("K132").Value = 2574.17652
$ws.Range("L132").Value = 15011.5005
$ws.Range("M132").Value = -44.17651999999998
$ws.Range("N132").Value = -20071.5005
$ws.Range("N39").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 583.3333
$ws.Range("I3").Value = 583.3333
$ws.Range("K3").Value = 583.3333
$ws.Range("M3").Value = -469.3333
$ws.Range("H80").Value = 376.33334
$ws.Range("I80").Value = 65
$ws.Range("J80").Value = 465.2857
$ws.Range("K80").Value = 65
$ws.Range("L80").Value = 465.2857
$ws.Range("M80").Value = 933
$ws.Range("N80").Value = -2461.2857
$ws.Range("H83").Value = 376.33334
$ws.Range("I83").Value = 65
$ws.Range("J83").Value = 465.2857
$ws.Range("K83").Value = 325
$ws.Range("L83").Value = 2326.4285
$ws.Range("M83").Value = 4667
$ws.Range("N83").Value = -12310.4285
$ws.Range("H99").Value = 3038.5
$ws.Range("I99").Value = 949.7778
$ws.Range("J99").Value = 4747.4546
$ws.Range("K99").Value = 949.7778
$ws.Range("L99").Value = 4747.4546
$ws.Range("M99").Value = 548.2222
$ws.Range("N99").Value = -7743.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3560.55
$ws.Range("I31").Value = 1169
$ws.Range("K31").Value = 1169
$ws.Range("M31").Value = -874
$ws.Range("H34").Value = 3560.55
$ws.Range("I34").Value = 1169
$ws.Range("K34").Value = 1169
$ws.Range("M34").Value = -967
$ws.Range("H58").Value = 2654.8406
$ws.Range("I58").Value = 1612.4364
$ws.Range("K58").Value = 1612.4364
$ws.Range("M58").Value = -1409.4364
$ws.Range("H99").Value = 12505000
$ws.Range("I99").Value = 33335600
$ws.Range("J99").Value = 6640
$ws.Range("K99").Value = 33335600
$ws.Range("L99").Value = 6640
$ws.Range("M99").Value = -33334102
$ws.Range("N99").Value = -9636
$ws.Range("H105").Value = 1497.5416
$ws.Range("I105").Value = 1187.6666
$ws.Range("K105").Value = 1187.6666
$ws.Range("M105").Value = 559.3334
$ws.Range("H126").Value = 12505000
$ws.Range("I126").Value = 33335600
$ws.Range("J126").Value = 6640
$ws.Range("K126").Value = 100006800
$ws.Range("L126").Value = 19920
$ws.Range("M126").Value = -100004330
$ws.Range("N126").Value = -24860
$ws.Range("H132").Value = 2465
$ws.Range("I132").Value = 1479.2778
$ws.Range("J132").Value = 4999.7144
$ws.Range("K132").Value = 4437.8334
$ws.Range("L132").Value = 14999.1432
$ws.Range("M132").Value = -1907.8334
$ws.Range("N132").Value = -20059.1432
$ws.Range("H134").Value = 5758.68
$ws.Range("I134").Value = 6739.353
$ws.Range("J134").Value = 3674.75
$ws.Range("K134").Value = 20218.059
$ws.Range("L134").Value = 11024.25
$ws.Range("M134").Value = -17683.059
$ws.Range("N134").Value = -16094.25
$ws.Range("H136").Value = 2654.8406
$ws.Range("I136").Value = 1612.4364
$ws.Range("K136").Value = 4837.3092
$ws.Range("M136").Value = -2287.3092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 22764.857
$ws.Range("J34").Value = 24444
$ws.Range("L34").Value = 73332
$ws.Range("N34").Value = -73500
$ws.Range("H113").Value = 3572003
$ws.Range("I113").Value = 601
$ws.Range("J113").Value = 6579499.5
$ws.Range("K113").Value = 1803
$ws.Range("L113").Value = 19738498.5
$ws.Range("M113").Value = 367
$ws.Range("N113").Value = -19742838.5
$ws.Range("H131").Value = 746.26
$ws.Range("I131").Value = 327.0909
$ws.Range("J131").Value = 798.06744
$ws.Range("K131").Value = 981.2727
$ws.Range("L131").Value = 2394.20232
$ws.Range("M131").Value = 4058.7273
$ws.Range("N131").Value = -12474.20232
$ws.Range("H132").Value = 2026.4762
$ws.Range("I132").Value = 927
$ws.Range("J132").Value = 3235.9
$ws.Range("K132").Value = 8343
$ws.Range("L132").Value = 29123.1
$ws.Range("M132").Value = -5813
$ws.Range("N132").Value = -34183.10000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3719.1
$ws.Range("I102").Value = 2198.5
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 2198.5
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -576.5
$ws.Range("N102").Value = -9244
$ws.Range("H132").Value = 2966.4075
$ws.Range("I132").Value = 1718.0555
$ws.Range("J132").Value = 5463.1113
$ws.Range("K132").Value = 5154.166499999999
$ws.Range("L132").Value = 16389.3339
$ws.Range("M132").Value = -2624.166499999999
$ws.Range("N132").Value = -21449.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 31067.334
$ws.Range("I10").Value = 203
$ws.Range("K10").Value = 203
$ws.Range("H12").Value = 1324.75
$ws.Range("I12").Value = 266.66666
$ws.Range("J12").Value = 4499
$ws.Range("K12").Value = 266.66666
$ws.Range("L12").Value = 4499
$ws.Range("M12").Value = -96.66665999999998
$ws.Range("N12").Value = -4839
$ws.Range("H17").Value = 3958.4285
$ws.Range("J17").Value = 4741.8
$ws.Range("L17").Value = 4741.8
$ws.Range("N17").Value = -5081.8
$ws.Range("H122").Value = 3906.2222
$ws.Range("I122").Value = 2108.2307
$ws.Range("J122").Value = 8581
$ws.Range("K122").Value = 6324.6921
$ws.Range("L122").Value = 25743
$ws.Range("M122").Value = -3874.6921
$ws.Range("N122").Value = -30643
$ws.Range("H136").Value = 4947.273
$ws.Range("I136").Value = 1784
$ws.Range("K136").Value = 5352
$ws.Range("M136").Value = -2802
$ws.Range("M10").Value = -63

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 50000
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("H13").Value = 245
$ws.Range("I13").Value = 245
$ws.Range("K13").Value = 245
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("H54").Value = 17894.572
$ws.Range("J54").Value = 19210.334
$ws.Range("L54").Value = 19210.334
$ws.Range("N54").Value = -20250.334
$ws.Range("H81").Value = 3254.5454
$ws.Range("I81").Value = 1666.6666
$ws.Range("J81").Value = 5160
$ws.Range("K81").Value = 3333.3332
$ws.Range("L81").Value = 10320
$ws.Range("M81").Value = -2272.3332
$ws.Range("N81").Value = -12442
$ws.Range("H84").Value = 3254.5454
$ws.Range("I84").Value = 1666.6666
$ws.Range("J84").Value = 5160
$ws.Range("K84").Value = 16666.666
$ws.Range("L84").Value = 51600
$ws.Range("M84").Value = -11362.666
$ws.Range("N84").Value = -62208
$ws.Range("H107").Value = 723.65
$ws.Range("I107").Value = 696.0625
$ws.Range("J107").Value = 834
$ws.Range("K107").Value = 2088.1875
$ws.Range("L107").Value = 2502
$ws.Range("M107").Value = -168.1875
$ws.Range("N107").Value = -6342
$ws.Range("H132").Value = 7578338
$ws.Range("I132").Value = 1697.2222
$ws.Range("K132").Value = 5091.6666
$ws.Range("M132").Value = -2561.6666
$ws.Range("M13").Value = -105
$ws.Range("M10").ClearContents()
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()
